$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Fgf1", "Cspg4", "ECs", 2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 2, 0.6666666666666666, 2.403560666666667, 7.210682, 0.06725806189087251, 0.06725806189087251, 4.818563884480445, 43.367074960324, 0.008075665595009171, 0.008075665595009171),
  @("ECs", "Fgf1", "Cspg4", "FAPs", 2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 3, 1, 12.844987, 38.534961, 0.3594371228547255, 0.3594371228547256, 25.75112470144466, 231.760122313002, 0.04315756245424775, 0.04315756245424777),
  @("ECs", "Fgf1", "Cspg4", "sCs", 2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 3, 1, 20.48784733333333, 61.463542, 0.5733048152544019, 0.573304815254402, 41.07323047853822, 369.659074306844, 0.06883662481257684, 0.06883662481257685),
  @("FAPs", "Fgf1", "Cspg4", "ECs", 3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 2, 0.6666666666666666, 2.403560666666667, 7.210682, 0.06725806189087251, 0.06725806189087251, 22.61667664490333, 203.55008980413, 0.037904388492794, 0.037904388492794),
  @("FAPs", "Fgf1", "Cspg4", "FAPs", 3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 3, 1, 12.844987, 38.534961, 0.3594371228547255, 0.3594371228547256, 120.866896149485, 1087.802065345365, 0.2025667103747836, 0.2025667103747836),
  @("FAPs", "Fgf1", "Cspg4", "sCs", 3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 3, 1, 20.48784733333333, 61.463542, 0.5733048152544019, 0.573304815254402, 192.7835750993366, 1735.05217589403, 0.3230953707445649, 0.3230953707445648),
  @("sCs", "Fgf1", "Cspg4", "ECs", 3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 2, 0.6666666666666666, 2.403560666666667, 7.210682, 0.06725806189087251, 0.06725806189087251, 12.69609776770933, 114.264879909384, 0.02127800780306934, 0.02127800780306934),
  @("sCs", "Fgf1", "Cspg4", "FAPs", 3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 3, 1, 12.844987, 38.534961, 0.3594371228547255, 0.3594371228547256, 67.84984171134798, 610.6485754021319, 0.1137128500256942, 0.1137128500256942),
  @("sCs", "Fgf1", "Cspg4", "sCs", 3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 3, 1, 20.48784733333333, 61.463542, 0.5733048152544019, 0.573304815254402, 108.2209891355227, 973.988902219704, 0.1813728196972602, 0.1813728196972603),
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
  }
}